# Regenerate the handback report: refresh the "Correspond Handoff Datetime"
# and "Correspond Handback DateTime" timestamps on the zh-cn and de-de
# language sheets (row 2 - the 7b66201e-... file pair).

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-11 03:06:44"   # Correspond Handoff Datetime
$wsZhCn.Range("G2").Value = "2016-01-11 03:07:26"   # Correspond Handback DateTime

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-11 03:06:55"   # Correspond Handoff Datetime
$wsDeDe.Range("G2").Value = "2016-01-11 03:07:45"   # Correspond Handback DateTime
